$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: "ainda não divulgado" becomes the Cap.7 summary assignment text,
# and the cell gains wrap text (style 22 -> 24).
$ws.Range("B17").Value = "Fazer um resumo (min 3 pgs) do Cap. 7 (seções 7.1, 7.5)"
$ws.Range("B17").WrapText = $true

# Rows 19/21/23 previously all shared the placeholder "ainda não divulgado"
# text; they become the per-chapter list assignments, picking up the
# "ShrinkToFit explicitly set" style variant (style 22 -> 12).
$ws.Range("B19").Value = "Lista Cap 9"
$ws.Range("B19").ShrinkToFit = $false

$ws.Range("B21").Value = "Lista Cap 10"
$ws.Range("B21").ShrinkToFit = $false

$ws.Range("B23").Value = "Lista Cap 11"
$ws.Range("B23").ShrinkToFit = $false
